$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Rule Name" header cell above the Rule 1 row (shares the
# default cell style with A9, consistent with the target workbook).
$ws.Range("A8").Value = "Rule Name"

# Update the saved view/selection: scroll back to show column A and
# select cell A9 (previously the view was scrolled to D1 with I13 selected).
$ws.Range("A1").Select()
$ws.Range("A9").Select()
